$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Update weekly crime statistics table (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 133.333333333333
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 27
$ws.Range("H14").Value = -7.407407407407
$ws.Range("I14").Value = 371
$ws.Range("J14").Value = 427
$ws.Range("K14").Value = -13.114754098360
$ws.Range("L14").Value = -11.876484560570
$ws.Range("M14").Value = -22.058823529411
$ws.Range("N14").Value = -78.253223915592

# Row 15
$ws.Range("C15").Value = 28
$ws.Range("D15").Value = 27
$ws.Range("E15").Value = 3.703703703703
$ws.Range("F15").Value = 122
$ws.Range("G15").Value = 112
$ws.Range("H15").Value = 8.928571428571
$ws.Range("I15").Value = 1440
$ws.Range("J15").Value = 1305
$ws.Range("K15").Value = 10.344827586206
$ws.Range("L15").Value = 13.029827315541
$ws.Range("M15").Value = 21.724429416737
$ws.Range("N15").Value = -49.87817612252

# Row 16
$ws.Range("C16").Value = 389
$ws.Range("D16").Value = 358
$ws.Range("E16").Value = 8.659217877094
$ws.Range("F16").Value = 1430
$ws.Range("G16").Value = 1352
$ws.Range("H16").Value = 5.769230769230
$ws.Range("I16").Value = 15303
$ws.Range("J16").Value = 11689
$ws.Range("K16").Value = 30.917957053640
$ws.Range("L16").Value = 34.816315743106
$ws.Range("M16").Value = -9.245641086466
$ws.Range("N16").Value = -79.460714573323

# Row 17
$ws.Range("C17").Value = 447
$ws.Range("D17").Value = 463
$ws.Range("E17").Value = -3.455723542116
$ws.Range("F17").Value = 1998
$ws.Range("G17").Value = 1911
$ws.Range("H17").Value = 4.552590266875
$ws.Range("I17").Value = 22824
$ws.Range("J17").Value = 20063
$ws.Range("K17").Value = 13.761650799980
$ws.Range("L17").Value = 24.809974298682
$ws.Range("M17").Value = 52.302148672094
$ws.Range("N17").Value = -37.542073721369

# Row 18
$ws.Range("C18").Value = 281
$ws.Range("D18").Value = 291
$ws.Range("E18").Value = -3.436426116838
$ws.Range("F18").Value = 1186
$ws.Range("G18").Value = 1162
$ws.Range("H18").Value = 2.065404475043
$ws.Range("I18").Value = 13694
$ws.Range("J18").Value = 10727
$ws.Range("K18").Value = 27.659177775706
$ws.Range("L18").Value = 1.814126394052
$ws.Range("M18").Value = -15.666954058381
$ws.Range("N18").Value = -84.339691688394

# Row 19
$ws.Range("C19").Value = 911
$ws.Range("D19").Value = 973
$ws.Range("E19").Value = -6.372045220966
$ws.Range("F19").Value = 4228
$ws.Range("G19").Value = 3886
$ws.Range("H19").Value = 8.800823468862
$ws.Range("I19").Value = 44966
$ws.Range("J19").Value = 32916
$ws.Range("K19").Value = 36.608336371369
$ws.Range("L19").Value = 44.413398850242
$ws.Range("M19").Value = 37.120726984417
$ws.Range("N19").Value = -39.724668570126

# Row 20
$ws.Range("C20").Value = 268
$ws.Range("D20").Value = 229
$ws.Range("E20").Value = 17.030567685589
$ws.Range("F20").Value = 1107
$ws.Range("G20").Value = 917
$ws.Range("H20").Value = 20.719738276990
$ws.Range("I20").Value = 11796
$ws.Range("J20").Value = 8818
$ws.Range("K20").Value = 33.771830347017
$ws.Range("L20").Value = 51.600051407274
$ws.Range("M20").Value = 30.848585690515
$ws.Range("N20").Value = -87.902779202133

# Row 21
$ws.Range("C21").Value = 2331
$ws.Range("D21").Value = 2344
$ws.Range("E21").Value = -0.554607508532
$ws.Range("F21").Value = 10096
$ws.Range("G21").Value = 9367
$ws.Range("H21").Value = 7.782641187146
$ws.Range("I21").Value = 110394
$ws.Range("J21").Value = 85945
$ws.Range("K21").Value = 28.447262784338
$ws.Range("L21").Value = 31.890897360844
$ws.Range("M21").Value = 20.579336559151
$ws.Range("N21").Value = -70.575958932041

# Row 22
$ws.Range("C22").Value = 38
$ws.Range("D22").Value = 52
$ws.Range("E22").Value = -26.923076923076
$ws.Range("F22").Value = 204
$ws.Range("G22").Value = 180
$ws.Range("H22").Value = 13.333333333333
$ws.Range("I22").Value = 2019
$ws.Range("J22").Value = 1462
$ws.Range("K22").Value = 38.098495212038
$ws.Range("L22").Value = 28.762755102040
$ws.Range("M22").Value = 8.958445763626

# Row 23
$ws.Range("C23").Value = 98
$ws.Range("D23").Value = 92
$ws.Range("E23").Value = 6.521739130434
$ws.Range("F23").Value = 438
$ws.Range("G23").Value = 449
$ws.Range("H23").Value = -2.449888641425
$ws.Range("I23").Value = 5179
$ws.Range("J23").Value = 4774
$ws.Range("K23").Value = 8.483452031839
$ws.Range("L23").Value = 15.037760995113
$ws.Range("M23").Value = 41.271140207310

# Row 24
$ws.Range("C24").Value = 2099
$ws.Range("D24").Value = 2061
$ws.Range("E24").Value = 1.843765162542
$ws.Range("F24").Value = 8943
$ws.Range("G24").Value = 7994
$ws.Range("H24").Value = 11.871403552664
$ws.Range("I24").Value = 100447
$ws.Range("J24").Value = 73271
$ws.Range("K24").Value = 37.089708070041
$ws.Range("L24").Value = 42.179537991167
$ws.Range("M24").Value = 40.889262921663

# Row 25
$ws.Range("C25").Value = 795
$ws.Range("D25").Value = 812
$ws.Range("E25").Value = -2.093596059113
$ws.Range("F25").Value = 3330
$ws.Range("G25").Value = 3133
$ws.Range("H25").Value = 6.287902968400
$ws.Range("I25").Value = 36255
$ws.Range("J25").Value = 31387
$ws.Range("K25").Value = 15.509605887787
$ws.Range("L25").Value = 23.978387990288
$ws.Range("M25").Value = -10.346447736096

# Row 26
$ws.Range("C26").Value = 45
$ws.Range("D26").Value = 47
$ws.Range("E26").Value = -4.255319148936
$ws.Range("F26").Value = 191
$ws.Range("G26").Value = 198
$ws.Range("H26").Value = -3.535353535353
$ws.Range("I26").Value = 2302
$ws.Range("J26").Value = 2143
$ws.Range("K26").Value = 7.419505366308
$ws.Range("L26").Value = 15.157578789394

# Row 27
$ws.Range("C27").Value = 89
$ws.Range("D27").Value = 112
$ws.Range("E27").Value = -20.535714285714
$ws.Range("F27").Value = 449
$ws.Range("G27").Value = 414
$ws.Range("H27").Value = 8.454106280193
$ws.Range("I27").Value = 4565
$ws.Range("J27").Value = 4283
$ws.Range("K27").Value = 6.584169974317
$ws.Range("L27").Value = 35.019225081336

# Row 28
$ws.Range("C28").Value = 30
$ws.Range("D28").Value = 41
$ws.Range("E28").Value = -26.829268292682
$ws.Range("G28").Value = 130
$ws.Range("H28").Value = -29.230769230769
$ws.Range("I28").Value = 1422
$ws.Range("J28").Value = 1651
$ws.Range("K28").Value = -13.870381586917
$ws.Range("L28").Value = -14.182257091128
$ws.Range("M28").Value = -9.885931558935
$ws.Range("N28").Value = -72.960638904734

# Row 29
$ws.Range("C29").Value = 23
$ws.Range("D29").Value = 36
$ws.Range("E29").Value = -36.111111111111
$ws.Range("F29").Value = 76
$ws.Range("G29").Value = 112
$ws.Range("H29").Value = -32.142857142857
$ws.Range("I29").Value = 1169
$ws.Range("J29").Value = 1378
$ws.Range("K29").Value = -15.166908563135
$ws.Range("L29").Value = -13.343217197924
$ws.Range("M29").Value = -10.215053763440
$ws.Range("N29").Value = -75.238296970980

# Row 30
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = -23.076923076923
$ws.Range("F30").Value = 44
$ws.Range("G30").Value = 44
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 551
$ws.Range("J30").Value = 475
$ws.Range("K30").Value = 16
$ws.Range("L30").Value = 131.512605042017
